$wb = $excel.ActiveWorkbook

# --- Add new sheet "MiniEPGScreen" after the last sheet (DTVChannel) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "MiniEPGScreen"

# Populate cells bottom-up / right-to-left so new shared-string entries are
# allocated in the same order as the target workbook.
$newSheet.Range("B4").Value = "volgende"
$newSheet.Range("A4").Value = "Future"
$newSheet.Range("B3").Value = "bezig"
$newSheet.Range("A3").Value = "InProgress"
$newSheet.Range("B2").Value = "vorige"
$newSheet.Range("A2").Value = "Previous"
$newSheet.Range("A1").Value = "objectID"
$newSheet.Range("B1").Value = "name_nl"

# Copy the bold header formatting (style used for row-1 headers elsewhere
# in the workbook) onto the new sheet's header row.
$screenTitles = $wb.Worksheets.Item("screenTitles")
$screenTitles.Range("A1:B1").Copy() | Out-Null
$newSheet.Range("A1:B1").PasteSpecial(-4122) | Out-Null

$newSheet.Range("N25").Select() | Out-Null

# --- Update screenTitles sheet: add a new "Search"/"zoeken" row ---
$screenTitles.Range("A15").Value = "Search"
$screenTitles.Range("B15").Value = "zoeken"
$screenTitles.Range("A14:B14").Copy() | Out-Null
$screenTitles.Range("A15:B15").PasteSpecial(-4122) | Out-Null
$screenTitles.Range("A15:B15").Select() | Out-Null

# --- Make the new sheet the active / selected tab ---
$newSheet.Activate()
